# Add files via upload
# Fills in the "Journal de travail" sheet with the latest work-log entries
# (rows 35-39, both the left block A:C and the right block E:G) and updates
# the running totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")
$ws.Activate()

# --- Left block (A:C) -------------------------------------------------
# Row 37 gains a new entry; rows 38/39 stay empty as before.
$ws.Range("A37").Value = 45100
$ws.Range("B37").Value = "Implémentation du client"
$ws.Range("C37").Value = 7

# --- Right block (E:G) -------------------------------------------------
$ws.Range("E35").Value = 45099
$ws.Range("F35").Value = "Revisions, préparations des mises en place"
$ws.Range("G35").Value = 1

$ws.Range("E36").Value = 45099
$ws.Range("F36").Value = "Implémentation client1"
$ws.Range("G36").Value = 2.5

$ws.Range("E37").Value = 45099
$ws.Range("F37").Value = "Documentation des TT"
$ws.Range("G37").Value = 3

$ws.Range("E38").Value = 45100
$ws.Range("F38").Value = "Implémentation client1"
$ws.Range("G38").Value = 5

$ws.Range("E39").Value = 45100
$ws.Range("F39").Value = "Implémentation API-Gateway"
$ws.Range("G39").Value = 2.5

# Move the view to roughly where the author left it after these edits.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F40").Select()
